$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 108, pushing existing rows 108..118 down to 109..119
$ws.Rows.Item(108).Insert()

# Fill in the new row 108 with a new weekly price observation
$ws.Range("A108").Value = 7
$ws.Range("B108").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C108").Value = "Ñuble"
$ws.Range("D108").Value = 44449
$ws.Range("E108").Value = 16
$ws.Range("F108").Value = 100112003
$ws.Range("G108").Value = "Ajo"
$ws.Range("H108").Value = "Chino"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 60
$ws.Range("K108").Value = 15000
$ws.Range("L108").Value = 16000
$ws.Range("M108").Value = 15500
$ws.Range("N108").Value = "$/caja 10 kilos"
$ws.Range("O108").Value = "China"
$ws.Range("P108").Value = 1550
$ws.Range("Q108").Value = 10
$ws.Range("R108").Value = "Hortaliza"
